$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended after the run on 2025-10-11.
# Force the date cell to stay plain text (matching the rest of column A)
# instead of being auto-parsed into a date serial by the "looks like a
# date" heuristic, then drop the temporary text format so the cell is
# left with no explicit style, matching the existing rows.
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "10/11/2025"
$ws.Range("A55").ClearFormats()

$ws.Range("B55").Value = 11071.3
